$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (Item 18)
$ws.Range("B15").Value = "open"
$ws.Range("C15").Value = 42270
$ws.Range("C15").NumberFormat = "d-mmm"
$ws.Range("D15").Value = 42263
$ws.Range("D15").NumberFormat = "d-mmm"
$ws.Range("E15").Value = "Bryce"
$ws.Range("F15").Value = "HIGH"
$ws.Range("G15").Value = "SIL"
$ws.Range("H15").Value = "SIL "

# Row 16 (Item 19)
$ws.Range("B16").Value = "open "
$ws.Range("C16").Value = 42270
$ws.Range("C16").NumberFormat = "d-mmm"
$ws.Range("D16").Value = 42263
$ws.Range("D16").NumberFormat = "d-mmm"
$ws.Range("E16").Value = "Prashant"
$ws.Range("F16").Value = "Medium"
$ws.Range("G16").Value = "Understand code currently onboard Pixhawk"

# Row 17 (Item 20)
$ws.Range("B17").Value = "open"
$ws.Range("C17").Value = 42270
$ws.Range("C17").NumberFormat = "d-mmm"
$ws.Range("D17").Value = 42263
$ws.Range("D17").NumberFormat = "d-mmm"
$ws.Range("E17").Value = "Matt  "
$ws.Range("F17").Value = "Medium"
$ws.Range("G17").Value = "Set up ROS and communication on Odroid"

# Row 18 (Item 21)
$ws.Range("A18").Value = 21
$ws.Range("B18").Value = "open "
$ws.Range("C18").Value = 42270
$ws.Range("C18").NumberFormat = "d-mmm"
$ws.Range("D18").Value = 42263
$ws.Range("D18").NumberFormat = "d-mmm"
$ws.Range("E18").Value = "Ed"
$ws.Range("F18").Value = "HIIGH"
$ws.Range("G18").Value = "Develop control block diagram"

# Selection change
$ws.Range("D15").Select() | Out-Null
